$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSuite")

# Set Runmode to "Y" for the Customer and Products suites so the whole suite runs.
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"
